$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 325, shifting existing rows 325-362 down to 326-363
$ws.Rows.Item(325).Insert()

# Populate the new row 325 with the new weekly record
$ws.Cells.Item(325, 1).Value = 7
$ws.Cells.Item(325, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(325, 3).Value = "Ñuble"
$ws.Cells.Item(325, 4).Value = 44918
$ws.Cells.Item(325, 5).Value = 16
$ws.Cells.Item(325, 6).Value = 100112008
$ws.Cells.Item(325, 7).Value = "Coliflor"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 500
$ws.Cells.Item(325, 11).Value = 800
$ws.Cells.Item(325, 12).Value = 900
$ws.Cells.Item(325, 13).Value = 850
$ws.Cells.Item(325, 14).Value = "$/unidad"
$ws.Cells.Item(325, 15).Value = "Región del Maule"
$ws.Cells.Item(325, 16).Value = 850
$ws.Cells.Item(325, 17).Value = 1
$ws.Cells.Item(325, 18).Value = "Hortaliza"
